$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4cb5c74bcbd7e23cf9deba24f7a8cf935d99e36c/e2e/ca6eecfb-577d-4341-b61f-f8ac843d3b80.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80b7bc2f6e1423d3abc081ee92a27e6af653d149/e2e/ca6eecfb-577d-4341-b61f-f8ac843d3b80.md."

# --- Overview sheet: row for ca6eecfb-577d-4341-b61f-f8ac843d3b80.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = "2016-09-02 00:59:18"

# --- zh-cn sheet: row for ca6eecfb-577d-4341-b61f-f8ac843d3b80.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-09-02 00:59:14"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row for ca6eecfb-577d-4341-b61f-f8ac843d3b80.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-09-02 00:59:18"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667

Write-Output "Applied handoff status update"
